$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.685.52'
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").Value = '''1.854.39'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '''264.83'
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = '''0.5238'
$ws.Range("E7").Value = '  +1.25%  '
$ws.Range("D8").Value = '''0.3277'
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("D9").Value = '''0.06795'
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").Value = '''18.85'
$ws.Range("E10").Value = '  -1.79%  '
$ws.Range("D11").Value = '''0.7764'
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("D12").Value = '''0.07767'
$ws.Range("E12").Value = '  +1.29%  '
$ws.Range("D13").Value = '''1.868.78'
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("D14").Value = '''88.63'
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("D15").Value = '''5.037'
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("D16").Value = '''1.001'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '''14.05'
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").Value = '''0.000007977'
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").Value = '''26.719.91'
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("D21").Value = '''2.095.29'
$ws.Range("E21").Value = '  +1.92%  '
$ws.Range("D22").Value = '''4.645'
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("D23").Value = '''9.537'
$ws.Range("E23").Value = '  -0.74%  '
$ws.Range("D24").Value = '''6.009'
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").Value = '''143.46'
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("D26").Value = '''2.205'
$ws.Range("E26").Value = '  -5.84%  '
$ws.Range("D27").Value = '''1.682'
$ws.Range("E27").Value = '  +1.71%  '
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("D29").Value = '''112.31'
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("D30").Value = '''4.203'
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("D31").Value = '''4.160'
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("D32").Value = '''0.08771'
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("D33").Value = '''0.04836'
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").Value = '''1.138'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '''0.7191'
$ws.Range("E35").Value = '  +4.70%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '''2.876'
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("D38").Value = '''0.01788'
$ws.Range("E38").Value = '  -0.90%  '
$ws.Range("D39").Value = '''2.207'
$ws.Range("E39").Value = '  -1.21%  '
$ws.Range("D40").Value = '''0.4900'
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D41").Value = '''112.97'
$ws.Range("E41").Value = '  +0.65%  '
$ws.Range("D42").Value = '''0.8978'
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("D43").Value = '''6.083'
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").Value = '''0.9997'
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").Value = '''7.718'
$ws.Range("E45").Value = '  -1.03%  '
$ws.Range("D46").Value = '''0.4191'
$ws.Range("E46").Value = '  -1.64%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.05930'
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''9.103'
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("D49").Value = '''0.1239'
$ws.Range("E49").Value = '  -3.68%  '
$ws.Range("D50").Value = '''35.10'
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").Value = '''0.8882'
$ws.Range("E51").Value = '  +3.20%  '
